$wb = $excel.ActiveWorkbook

$schemaSheet = $wb.Worksheets.Item("!!_Schema")
$compoundSheet = $wb.Worksheets.Item("!!Compound")
$modelSheet = $wb.Worksheets.Item("!!Model")
$reactionSheet = $wb.Worksheets.Item("!!Reaction")

$schemaSheet.Unprotect()
$compoundSheet.Unprotect()
$modelSheet.Unprotect()
$reactionSheet.Unprotect()

$schemaSheet.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8' date='2020-03-09 15:30:10'"
$schemaSheet.Range("A2").Value = "!!ObjTables type='Schema' tableFormat='row' description='Table/model and column/attribute definitions' date='2020-03-09 15:30:10' objTablesVersion='0.0.8'"

$compoundSheet.Range("A1").Value = "!!ObjTables type='Data' id='Compound' description='Compound' name='Compound' date='2020-03-09 15:30:10' objTablesVersion='0.0.8' tableFormat='row'"

$modelSheet.Range("A1").Value = "!!ObjTables type='Data' id='Model' description='Model' name='Model' date='2020-03-09 15:30:10' objTablesVersion='0.0.8' tableFormat='column'"

$reactionSheet.Range("A1").Value = "!!ObjTables type='Data' id='Reaction' description='Reaction' name='Reaction' date='2020-03-09 15:30:10' objTablesVersion='0.0.8' tableFormat='row'"

$schemaSheet.Protect($null, $true, $true, $true)
$compoundSheet.Protect($null, $true, $true, $true)
$modelSheet.Protect($null, $true, $true, $true)
$reactionSheet.Protect($null, $true, $true, $true)
